# Fine tune the "person-create" contact form: split Name into First/Last
# Name, add a Current Address field, tweak the "Also Known As:" label
# casing, and drop the old "digital enrollment age" question.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# 1) "Name" -> "First Name" (row 19 keeps its place).
$ws.Range("C19").Value = "First Name"

# 2) Insert a new row right below First Name for the new "Last Name"
#    question (pushes aka/docket/tel/... down by one row).
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "string"
$ws.Range("B20").Value = "name1"
$ws.Range("C20").Value = "Last Name"

# 3) The "aka" row is now row 21 - fix label casing: "As" -> "as".
$ws.Range("C21").Value = "Also Known as:"

# 4) Insert a new row after "Docket Number" (now row 22) for the new
#    "Current Address" question (pushes tel/phone/... down by one row).
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "string "
$ws.Range("B23").Value = "address"
$ws.Range("C23").Value = "Current Address"

# 5) Remove the old "digital" / enrollment-age question, now at row 28
#    (sex/gender/status/meta group shift back up by one row).
$ws.Rows.Item(28).Delete()

# 6) Add a second "end group" row at the very end of the sheet to close
#    the newly balanced group structure.
$ws.Range("A36").Value = "end group"
